# The edit rotates the data of rows 2-5 (the 4 observation records) by one
# position: new row2 <- old row5, new row3 <- old row4, new row4 <- old row3,
# new row5 <- old row2. The Ost/Nord (Q/R) coordinate columns additionally
# get rounded to whole numbers as part of this update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually carry row-specific data in rows 2-5.
$cols = @("A","B","D","E","F","G","H","Q","R","S","Z","AB","AW","AX")

# Capture the "before" values for every relevant cell in rows 2-5 so we can
# rotate them without clobbering source data while we write.
$before = @{}
foreach ($r in 2..5) {
    foreach ($col in $cols) {
        $before["$col$r"] = $ws.Range("$col$r").Value2
    }
}

# Mapping of destination row -> source row (rotate up by one, wrapping 2<-5).
$rowMap = @{ 2 = 5; 3 = 4; 4 = 3; 5 = 2 }

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $val = $before["$col$srcRow"]
        if ($col -eq "Q" -or $col -eq "R") {
            $val = [Math]::Round([double]$val)
        }
        $ws.Range("$col$destRow").Value2 = $val
    }
}
